$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update branch names (column A) to reflect corrected/renamed branches.
$ws.Range("A12").Value = "Đặng Nguyên Cẩn"
$ws.Range("A16").Value = "3/26 Quang Trung"
$ws.Range("A23").Value = "Công Trường Quốc Tế"

# Reflect the user's last selection being on the edited row.
$ws.Range("A23").Select()
